$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.365.14"
$ws.Range("E2").Value = "  +0.16%  "
Set-TextValue "D3" "1.839.52"
$ws.Range("E3").Value = "  -0.16%  "
Set-TextValue "D4" "0.9992"
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue "D5" "238.93"
$ws.Range("E5").Value = "  -0.34%  "
Set-TextValue "D6" "0.6261"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.21%  "
Set-TextValue "D8" "0.07408"
$ws.Range("E8").Value = "  -0.85%  "
Set-TextValue "D9" "0.2890"
$ws.Range("E9").Value = "  -0.15%  "
Set-TextValue "D10" "24.93"
$ws.Range("E10").Value = "  +1.99%  "
Set-TextValue "D11" "0.07715"
$ws.Range("E11").Value = "  -0.23%  "
Set-TextValue "D12" "1.845.88"
$ws.Range("E12").Value = "  +0.18%  "
Set-TextValue "D13" "4.968"
$ws.Range("E13").Value = "  -0.18%  "
Set-TextValue "D14" "0.6735"
$ws.Range("E14").Value = "  -0.79%  "
Set-TextValue "D15" "0.00001024"
$ws.Range("E15").Value = "  -3.11%  "
Set-TextValue "D16" "81.71"
$ws.Range("E16").Value = "  -0.19%  "
Set-TextValue "D17" "6.202"
$ws.Range("E17").Value = "  +0.60%  "
Set-TextValue "D18" "29.410.21"
$ws.Range("E18").Value = "  +0.26%  "
Set-TextValue "D19" "232.99"
$ws.Range("E19").Value = "  +1.74%  "
Set-TextValue "D20" "12.31"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +0.17%  "
Set-TextValue "D22" "7.288"
$ws.Range("E22").Value = "  -2.74%  "
Set-TextValue "D23" "1.000"
$ws.Range("E23").Value = "  +0.10%  "
Set-TextValue "D24" "157.91"
$ws.Range("E24").Value = "  -0.29%  "
Set-TextValue "D25" "8.481"
$ws.Range("E25").Value = "  +0.75%  "
Set-TextValue "D26" "0.1344"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("E27").Value = "  -1.25%  "
Set-TextValue "D28" "0.07257"
$ws.Range("E28").Value = "  +10.28%  "
Set-TextValue "D29" "1.471"
$ws.Range("E29").Value = "  +4.14%  "
Set-TextValue "D30" "1.474"
$ws.Range("E30").Value = "  -0.53%  "
Set-TextValue "D31" "4.031"
$ws.Range("E31").Value = "  -1.83%  "
Set-TextValue "D32" "4.023"
$ws.Range("E32").Value = "  -1.53%  "
Set-TextValue "D33" "1.811"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  -0.03%  "
Set-TextValue "D35" "0.6953"
$ws.Range("E35").Value = "  +0.21%  "
Set-TextValue "D36" "2.571"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  +0.29%  "
Set-TextValue "D38" "6.901"
$ws.Range("E38").Value = "  +1.82%  "
Set-TextValue "D39" "2.813"
$ws.Range("E39").Value = "  -0.54%  "
Set-TextValue "D40" "1.233.21"
$ws.Range("E40").Value = "  -2.39%  "
Set-TextValue "D41" "0.9554"
$ws.Range("E41").Value = "  +4.07%  "
Set-TextValue "D42" "1.000"
$ws.Range("E42").Value = "  +0.24%  "
Set-TextValue "D43" "2.011.73"
$ws.Range("E43").Value = "  +0.42%  "
Set-TextValue "D44" "100.74"
$ws.Range("E44").Value = "  -0.45%  "
Set-TextValue "D45" "65.24"
$ws.Range("E45").Value = "  -1.24%  "
Set-TextValue "D48" "6.935"
$ws.Range("E48").Value = "  -1.81%  "
Set-TextValue "D49" "8.844"
$ws.Range("E49").Value = "  -1.28%  "
Set-TextValue "D50" "0.3894"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  -2.63%  "

# Rows 46/47 swap: RenderToken moves to row 46, BabyDogeCoin moves to row 47
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D46" "1.709"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.00000000116"
$ws.Range("E47").Value = "  -2.33%  "
